$d = $word.ActiveDocument

$rng = $d.Content
[void]$rng.Find.Execute("2023-07-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-21 Friday", 1)
$rng = $d.Content
[void]$rng.Find.Execute("39+40=", $true, $false, $false, $false, $false, $true, 1, $false, "72-34=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("11+37=", $true, $false, $false, $false, $false, $true, 1, $false, "82-53=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("78-39=", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("73+10=", $true, $false, $false, $false, $false, $true, 1, $false, "99-56=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("9+23=", $true, $false, $false, $false, $false, $true, 1, $false, "77+14=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("23-6=", $true, $false, $false, $false, $false, $true, 1, $false, "64-56=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("6+68=", $true, $false, $false, $false, $false, $true, 1, $false, "71-6=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("29+22=", $true, $false, $false, $false, $false, $true, 1, $false, "10+41=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("84-55=", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("94-17=", $true, $false, $false, $false, $false, $true, 1, $false, "70-58=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("42-20=", $true, $false, $false, $false, $false, $true, 1, $false, "74-3=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("94-60=", $true, $false, $false, $false, $false, $true, 1, $false, "13+85=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("15+29=", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("96-15=", $true, $false, $false, $false, $false, $true, 1, $false, "3+0=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("80-35=", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("49-42=", $true, $false, $false, $false, $false, $true, 1, $false, "8+5=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("36-0=", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("62+9=", $true, $false, $false, $false, $false, $true, 1, $false, "47-17=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("89-25=", $true, $false, $false, $false, $false, $true, 1, $false, "36-17=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("6+60=", $true, $false, $false, $false, $false, $true, 1, $false, "8+86=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("96-91=", $true, $false, $false, $false, $false, $true, 1, $false, "41-10=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("39-23=", $true, $false, $false, $false, $false, $true, 1, $false, "34+35=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("81-70=", $true, $false, $false, $false, $false, $true, 1, $false, "26-5=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("88-20=", $true, $false, $false, $false, $false, $true, 1, $false, "69-20=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("66-62=", $true, $false, $false, $false, $false, $true, 1, $false, "73+19=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("33-11=", $true, $false, $false, $false, $false, $true, 1, $false, "80-64=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("65+33=", $true, $false, $false, $false, $false, $true, 1, $false, "41+1=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("29+66=", $true, $false, $false, $false, $false, $true, 1, $false, "47+19=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("15+58=", $true, $false, $false, $false, $false, $true, 1, $false, "17+75=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("68-14=", $true, $false, $false, $false, $false, $true, 1, $false, "28-7=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("14+69=", $true, $false, $false, $false, $false, $true, 1, $false, "83-17=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("52+16=", $true, $false, $false, $false, $false, $true, 1, $false, "90-61=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("10+33=", $true, $false, $false, $false, $false, $true, 1, $false, "17+8=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("87-60=", $true, $false, $false, $false, $false, $true, 1, $false, "49-27=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("33+12=", $true, $false, $false, $false, $false, $true, 1, $false, "93-33=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("68+16=", $true, $false, $false, $false, $false, $true, 1, $false, "52-21=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("92-62=", $true, $false, $false, $false, $false, $true, 1, $false, "30+32=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("88-74=", $true, $false, $false, $false, $false, $true, 1, $false, "22+36=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("0+8=", $true, $false, $false, $false, $false, $true, 1, $false, "42+27=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("88-25=", $true, $false, $false, $false, $false, $true, 1, $false, "32+36=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("73-8=", $true, $false, $false, $false, $false, $true, 1, $false, "36+21=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("79-7=", $true, $false, $false, $false, $false, $true, 1, $false, "94-68=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("8+70=", $true, $false, $false, $false, $false, $true, 1, $false, "13+55=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("69+8=", $true, $false, $false, $false, $false, $true, 1, $false, "41-32=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("26-17=", $true, $false, $false, $false, $false, $true, 1, $false, "11+64=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("44+50=", $true, $false, $false, $false, $false, $true, 1, $false, "71+25=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("95-5=", $true, $false, $false, $false, $false, $true, 1, $false, "96-27=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("7+89=", $true, $false, $false, $false, $false, $true, 1, $false, "93-31=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("10+36=", $true, $false, $false, $false, $false, $true, 1, $false, "77-39=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("56+14=", $true, $false, $false, $false, $false, $true, 1, $false, "59+32=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("64+10=", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("56-13=", $true, $false, $false, $false, $false, $true, 1, $false, "19-19=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("72+19=", $true, $false, $false, $false, $false, $true, 1, $false, "82-49=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("4+64=", $true, $false, $false, $false, $false, $true, 1, $false, "36+43=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("73-42=", $true, $false, $false, $false, $false, $true, 1, $false, "5+40=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("27+3=", $true, $false, $false, $false, $false, $true, 1, $false, "6+80=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("46-15=", $true, $false, $false, $false, $false, $true, 1, $false, "23+14=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("55-25=", $true, $false, $false, $false, $false, $true, 1, $false, "89-71=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("19+65=", $true, $false, $false, $false, $false, $true, 1, $false, "10+61=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("96-70=", $true, $false, $false, $false, $false, $true, 1, $false, "21+71=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("73-43=", $true, $false, $false, $false, $false, $true, 1, $false, "17+28=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("89+8=", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("57+19=", $true, $false, $false, $false, $false, $true, 1, $false, "33+50=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("68-16=", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("23+1=", $true, $false, $false, $false, $false, $true, 1, $false, "33+19=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("15+65=", $true, $false, $false, $false, $false, $true, 1, $false, "37+16=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("72-38=", $true, $false, $false, $false, $false, $true, 1, $false, "78-25=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("98-79=", $true, $false, $false, $false, $false, $true, 1, $false, "50-21=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("49-44=", $true, $false, $false, $false, $false, $true, 1, $false, "53+36=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("56+1=", $true, $false, $false, $false, $false, $true, 1, $false, "85-44=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("27+14=", $true, $false, $false, $false, $false, $true, 1, $false, "5+47=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("33-25=", $true, $false, $false, $false, $false, $true, 1, $false, "47-12=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("93-71=", $true, $false, $false, $false, $false, $true, 1, $false, "43+50=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("67-51=", $true, $false, $false, $false, $false, $true, 1, $false, "62+8=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("88-13=", $true, $false, $false, $false, $false, $true, 1, $false, "97-41=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("80-20=", $true, $false, $false, $false, $false, $true, 1, $false, "14+71=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("57+19=", $true, $false, $false, $false, $false, $true, 1, $false, "3+28=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("52-37=", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("8+80=", $true, $false, $false, $false, $false, $true, 1, $false, "54-48=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("23-12=", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("36-5=", $true, $false, $false, $false, $false, $true, 1, $false, "20+39=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("87-40=", $true, $false, $false, $false, $false, $true, 1, $false, "23+69=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("18+4=", $true, $false, $false, $false, $false, $true, 1, $false, "98-19=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("87-4=", $true, $false, $false, $false, $false, $true, 1, $false, "90-8=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("49+42=", $true, $false, $false, $false, $false, $true, 1, $false, "1+28=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("68-68=", $true, $false, $false, $false, $false, $true, 1, $false, "23+44=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("29+36=", $true, $false, $false, $false, $false, $true, 1, $false, "99-74=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("11+68=", $true, $false, $false, $false, $false, $true, 1, $false, "14-0=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("18+79=", $true, $false, $false, $false, $false, $true, 1, $false, "58+36=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("54+40=", $true, $false, $false, $false, $false, $true, 1, $false, "48+30=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("80-33=", $true, $false, $false, $false, $false, $true, 1, $false, "86-78=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("27+54=", $true, $false, $false, $false, $false, $true, 1, $false, "78+5=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("48+25=", $true, $false, $false, $false, $false, $true, 1, $false, "87-70=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("2+29=", $true, $false, $false, $false, $false, $true, 1, $false, "73+19=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("96-52=", $true, $false, $false, $false, $false, $true, 1, $false, "57-56=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("47-15=", $true, $false, $false, $false, $false, $true, 1, $false, "41+0=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("26+7=", $true, $false, $false, $false, $false, $true, 1, $false, "24-1=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("89-10=", $true, $false, $false, $false, $false, $true, 1, $false, "41+53=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("60-38=", $true, $false, $false, $false, $false, $true, 1, $false, "1+35=", 1)
$rng = $d.Content
[void]$rng.Find.Execute("95-3=", $true, $false, $false, $false, $false, $true, 1, $false, "19+41=", 1)
